# "Bypass empty keyword cell"
#
# Adds a third Robot-Framework-style test suite sheet ("suite3") that
# mirrors "suite2" (Checking available flights) but merges the repeated
# "open flights page" keyword cell down the rows it spans, instead of
# repeating it (or leaving it blank) on every row.

$wb = $excel.ActiveWorkbook

# --- suite2: selection moved off the single C8 cell onto A1:D7 ---------
$ws2 = $wb.Worksheets.Item("suite2")
$ws2.Range("A1:D7").Select()

# --- suite3: new sheet appended after suite2 ----------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "suite3"

$ws3.Range("A1").Value = "test case"
$ws3.Range("B1").Value = "steps"
$ws3.Range("C1").Value = "arg1"
$ws3.Range("D1").Value = "arg2"

$ws3.Range("A2").Value = "Checking available flights"

# Keyword cell spans rows 3-6 (merged) instead of repeating/being blank
$ws3.Range("B3:B6").VerticalAlignment = -4160   # xlTop
$ws3.Range("B3:B6").Merge()
$ws3.Range("B3").Value = "open flights page"

$ws3.Range("B7").Value = "check flights to"
$ws3.Range("C7").Value = "Dublin"
$ws3.Range("B8").Value = "check flights to"
$ws3.Range("C8").Value = "New York"

$ws3.Columns.Item(1).AutoFit()
$ws3.Columns.Item(2).AutoFit()
$ws3.Columns.Item(3).AutoFit()

$ws3.Range("G12").Select()
$ws3.Activate()
